$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update values
$ws.Range("B4").Value = 1.98
$ws.Range("B5").Value = 2.3

# Update selection to B7
$ws.Range("B7").Select()
